$d = $word.ActiveDocument

# Helper: insert a new bulleted paragraph (same list formatting as the
# paragraph ending at $afterRange) right after $afterRange, containing $text.
# Returns the Range of the newly inserted paragraph's text.
function Insert-BulletAfter($afterRange, $text) {
    $pEnd = $afterRange.End
    $insertRng = $d.Range($pEnd - 1, $pEnd - 1)
    $insertRng.InsertParagraphAfter()
    $newParaRng = $d.Range($pEnd, $pEnd)
    $newPara = $newParaRng.Paragraphs(1)
    $newPara.Range.InsertAfter($text)
    return $newPara.Range
}

# 1) Update the "Mapping : ..." bullet text and add three new bullets
#    (SubjectMapping / PredicateMapping / ObjectMapping) right after it.
$rng = $d.Content
$found = $rng.Find.Execute("Mapping : (URN, Occurrence, AttributeKind, Value) : Statement;", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Mapping : (URN, Occurrence, Kind, Value) : Statement;", 2)

$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Mapping : (URN, Occurrence, Kind, Value) : Statement;")
$mappingPara = $rng2.Paragraphs(1)

$r1 = Insert-BulletAfter $mappingPara.Range "SubjectMapping : (URN, Occurrence, SubjectKind, Value) : Mapping;"
$r2 = Insert-BulletAfter $r1 "PredicateMapping : (URN, Occurrence, PredicateKind, Value) : Mapping;"
$r3 = Insert-BulletAfter $r2 "ObjectMapping : (URN, Occurrence, ObjectKind, Value) : Mapping;"

# 2) Append ". Regression." to the Activation bullet line.
$d.Content.Find.Execute(
    "Activation : Occurrences / Mappings. (Value in Context. Prediction: Speed, Time : Distance)", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Activation : Occurrences / Mappings. (Value in Context. Prediction: Speed, Time : Distance). Regression.", 2)
